$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -0.4267425320056861
$ws.Range("G3").Value = -0.421744981281047
$ws.Range("H3").Value = 98.82890728019288
$ws.Range("F4").Value = -0.4249291784702569
$ws.Range("G4").Value = -0.4100613845468226
$ws.Range("H4").Value = 96.50111249668515
$ws.Range("F5").Value = -0.4231311706629048
$ws.Range("G5").Value = -0.3985755992018447
$ws.Range("H5").Value = 94.19669994470277
$ws.Range("F6").Value = -0.4213483146067398
$ws.Range("G6").Value = -0.3872809822789258
$ws.Range("H6").Value = 91.91468646086545
$ws.Range("F7").Value = -0.4195804195804231
$ws.Range("G7").Value = -0.3761711465014916
$ws.Range("H7").Value = 89.6541232495214
$ws.Range("F8").Value = -0.417827298050133
$ws.Range("G8").Value = -0.3652399474717738
$ws.Range("H8").Value = 87.41409409491251
$ws.Range("F9").Value = -0.4160887656033285
$ws.Range("G9").Value = -0.3544814716036737
$ws.Range("H9").Value = 85.19371367541628
$ws.Range("F10").Value = -0.4143646408839907
$ws.Range("G10").Value = -0.3438900247496601
$ws.Range("H10").Value = 82.9921259729154
$ws.Range("F11").Value = -0.4126547455295726
$ws.Range("G11").Value = -0.3334601214759148
$ws.Range("H11").Value = 80.80850277099687
$ws.Range("F12").Value = -0.4109589041095929
$ws.Range("G12").Value = -0.3231864749417612
$ws.Range("H12").Value = 78.64204223582783
$ws.Range("F13").Value = -0.4092769440654842
$ws.Range("G13").Value = -0.3130639873435648
$ws.Range("H13").Value = 76.49196757427768
$ws.Range("F14").Value = -0.4076086956521729
$ws.Range("G14").Value = -0.3030877408869204
$ws.Range("H14").Value = 74.35752576425797
$ws.Range("F15").Value = -0.4059539918809141
$ws.Range("G15").Value = -0.2932529892521019
$ws.Range("H15").Value = 72.23798635243551
$ws.Range("F16").Value = -0.4043126684636134
$ws.Range("G16").Value = -0.2835551495213418
$ws.Range("H16").Value = 70.13264031494494
$ws.Range("F17").Value = -0.4026845637583931
$ws.Range("G17").Value = -0.2739897945386094
$ws.Range("H17").Value = 68.04079897708735
$ws.Range("F18").Value = -0.4010695187165791
$ws.Range("G18").Value = -0.2645526456748981
$ws.Range("H18").Value = 65.96179298827434
$ws.Range("F19").Value = -0.3994673768308865
$ws.Range("G19").Value = -0.2552395659724649
$ws.Range("H19").Value = 63.89497134844127
$ws.Range("F20").Value = -0.3978779840848823
$ws.Range("G20").Value = -0.2460465536460621
$ws.Range("H20").Value = 61.83970048304334
$ws.Range("F21").Value = -0.3963011889035695
$ws.Range("G21").Value = -0.2369697359177358
$ws.Range("H21").Value = 59.79536336324159
$ws.Range("F22").Value = -0.3947368421052588
$ws.Range("G22").Value = -0.2280053631656176
$ws.Range("H22").Value = 57.76135866862375
$ws.Range("F23").Value = -0.3931847968545177
$ws.Range("G23").Value = -0.2191498033675021
$ws.Range("H23").Value = 55.73709998980193
$ws.Range("F24").Value = -0.3916449086161844
$ws.Range("G24").Value = -0.2103995368209155
$ws.Range("H24").Value = 53.72201506827425
$ws.Range("F25").Value = -0.3901170351105376
$ws.Range("G25").Value = -0.2017511511241521
$ws.Range("H25").Value = 51.7155450714904
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
